$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "27.902.98"
$ws.Range("E2").Value = "  -0.05%  "
$ws.Range("D3").Value = "1.623.20"
$ws.Range("E3").Value = "  -0.39%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'211.37"
$ws.Range("E5").Value = "  -0.25%  "
$ws.Range("E6").Value = "  -1.48%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").Value = "'23.01"
$ws.Range("E8").Value = "  -1.54%  "
$ws.Range("E9").Value = "  +0.33%  "
$ws.Range("E10").Value = "  -1.31%  "
$ws.Range("E11").Value = "  +0.08%  "
$ws.Range("D12").Value = "1.852.92"
$ws.Range("E12").Value = "  -0.42%  "
$ws.Range("D13").Value = "1.627.98"
$ws.Range("E13").Value = "  +0.07%  "
$ws.Range("E14").Value = "  -0.84%  "
$ws.Range("D15").Value = "'0.552"
$ws.Range("D16").Value = "'64.49"
$ws.Range("E16").Value = "  -1.73%  "
$ws.Range("D17").Value = "27.890.51"
$ws.Range("E17").Value = "  -0.09%  "
$ws.Range("D18").Value = "'227.74"
$ws.Range("E18").Value = "  -1.14%  "
$ws.Range("D19").Value = "'7.59"
$ws.Range("E19").Value = "  -0.54%  "
$ws.Range("D20").Value = "0.0₃0715"
$ws.Range("E20").Value = "  -1.13%  "
$ws.Range("E21").Value = "  +0.20%  "
$ws.Range("E22").Value = "  -0.36%  "
$ws.Range("E23").Value = "  -3.24%  "
$ws.Range("D24").Value = "'2.07"
$ws.Range("E24").Value = "  +1.60%  "
$ws.Range("D25").Value = "'154.27"
$ws.Range("E25").Value = "  -0.37%  "
$ws.Range("D26").Value = "'6.90"
$ws.Range("E26").Value = "  -0.43%  "
$ws.Range("E27").Value = "  -0.61%  "
$ws.Range("E28").Value = "  +0.13%  "
$ws.Range("D29").Value = "'15.37"
$ws.Range("E29").Value = "  -1.12%  "
$ws.Range("E30").Value = "  -0.23%  "
$ws.Range("E31").Value = "  -0.58%  "
$ws.Range("E32").Value = "  -0.64%  "
$ws.Range("D33").Value = "1.415.09"
$ws.Range("E33").Value = "  +1.00%  "
$ws.Range("E34").Value = "  +0.17%  "
$ws.Range("E35").Value = "  +2.04%  "
$ws.Range("D36").Value = "'0.989"
$ws.Range("E36").Value = "  -2.57%  "
$ws.Range("E37").Value = "  -0.53%  "
$ws.Range("E38").Value = "  -0.85%  "
$ws.Range("D39").Value = "'0.554"
$ws.Range("E39").Value = "  -0.31%  "
$ws.Range("D40").Value = "'0.847"
$ws.Range("E40").Value = "  -1.84%  "
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("D42").Value = "'1.01"
$ws.Range("E42").Value = "  -1.76%  "
$ws.Range("D43").Value = "'65.28"
$ws.Range("E43").Value = "  -1.59%  "
$ws.Range("D44").Value = "'1.80"
$ws.Range("E44").Value = "  -2.34%  "
$ws.Range("D45").Value = "'5.36"
$ws.Range("E45").Value = "  -2.24%  "
$ws.Range("D46").Value = "1.762.98"
$ws.Range("E46").Value = "  -0.44%  "
$ws.Range("E47").Value = "  -3.88%  "
$ws.Range("D48").Value = "'89.10"
$ws.Range("E48").Value = "  +1.23%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₆0102"
$ws.Range("E49").Value = "  +0.89%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "'0.100"
$ws.Range("E50").Value = "  -0.42%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.0502"
$ws.Range("E51").Value = "  -0.49%  "
